$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.805.62'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '2.608.19'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('D9').Value = '2.633.28'
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.106'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.157'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.87%  '
$ws.Range('E13').Value = '  +5.46%  '
$ws.Range('D14').Value = '3.071.51'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '60.861.11'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('E17').Value = '  +4.19%  '
$ws.Range('D18').Value = '2.620.85'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '350.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.518'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.07%  '
$ws.Range('D29').Value = '0.0₃0801'
$ws.Range('E29').Value = '  +3.39%  '
$ws.Range('E30').Value = '  +9.10%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '162.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.02'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.72%  '
$ws.Range('E36').Value = '  +4.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.24'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.62'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.90'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('E40').Value = '  +5.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '306.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.847'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '133.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.36%  '
$ws.Range('E46').Value = '  +11.84%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0986'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('E49').Value = '  +2.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0551'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.62%  '
$ws.Range('E51').Value = '  +3.51%  '
